# GA_training_scenarios.xlsx - "aspiration example set up"
#
# 1. Rename two header cells on the existing "Tabelle1" sheet:
#      A1: id_training_scenario -> id
#      B1: number_of_training_path -> number_of_path
# 2. Add a new worksheet "Sheet1" right after "Tabelle1" that holds the
#    per-path training data (id_scenario / id_training_scenario / id_path /
#    id_agent / strategy_param_1..3).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Tabelle1")

# --- Tabelle1: rename the two header cells -------------------------------
$ws1.Range("A1").Value = "id"
$ws1.Range("B1").Value = "number_of_path"

# --- Add the new sheet, positioned right after Tabelle1 ------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet1"

# --- New sheet header row -------------------------------------------------
# (cell order below matches the order the string values were first written
#  to the workbook's shared-string table in the source workbook)
$ws2.Range("A1").Value = "id_scenario"
$ws2.Range("D1").Value = "id_agent"
$ws2.Range("E1").Value = "strategy_param_1"
$ws2.Range("F1").Value = "strategy_param_2"
$ws2.Range("G1").Value = "strategy_param_3"
$ws2.Range("C1").Value = "id_path"
$ws2.Range("B1").Value = "id_training_scenario"

# --- Restore the view state: selection on the new sheet, then reactivate -
# --- Tabelle1 so it remains the visible / selected tab -------------------
[void]$ws2.Range("E19").Select()
[void]$ws1.Activate()
[void]$ws1.Range("C16").Select()
